# Automatic update of files.
# Bump the "Förändrad" (Changed) date column (C) from 45745 to 45746
# for every data row (rows 2 through 43) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45745) {
        $cell.Value = 45746
    }
}
